$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column N (last existing year) into new column O, then set the new 2021 values.
$ws.Range("N3:N5").Copy()
$ws.Range("O3:O5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("O3").Value = 2021
$ws.Range("O4").Value = 14
$ws.Range("O5").Value = 1252.8

# Update the active selection to match the target workbook state.
$ws.Range("O9").Select()
